$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Sheet1 -> ValidLogin
$ws.Name = "ValidLogin"

# Populate the login form data
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Match the saved view state: active cell C2 selected, zoomed to 235%
[void]$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 235
